$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2023-09-26 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-09-27 Wednesday", 2)

# Update the division-problem table. Only rows 1, 5, 9, 13, 17 (1-based)
# carry answer text; the rest are blank spacer rows.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "40÷6=6, 4"
$t.Cell(1, 2).Range.Text = "38÷3=12, 2"
$t.Cell(1, 3).Range.Text = "98÷9=10, 8"
$t.Cell(1, 4).Range.Text = "68÷3=22, 2"
$t.Cell(1, 5).Range.Text = "97÷8=12, 1"

$t.Cell(5, 1).Range.Text = "15÷3=5, 0"
$t.Cell(5, 2).Range.Text = "65÷5=13, 0"
$t.Cell(5, 3).Range.Text = "21÷9=2, 3"
$t.Cell(5, 4).Range.Text = "23÷3=7, 2"
$t.Cell(5, 5).Range.Text = "54÷9=6, 0"

$t.Cell(9, 1).Range.Text = "74÷3=24, 2"
$t.Cell(9, 2).Range.Text = "98÷2=49, 0"
$t.Cell(9, 3).Range.Text = "31÷2=15, 1"
$t.Cell(9, 4).Range.Text = "94÷6=15, 4"
$t.Cell(9, 5).Range.Text = "33÷2=16, 1"

$t.Cell(13, 1).Range.Text = "50÷2=25, 0"
$t.Cell(13, 2).Range.Text = "29÷5=5, 4"
$t.Cell(13, 3).Range.Text = "25÷2=12, 1"
$t.Cell(13, 4).Range.Text = "64÷9=7, 1"
$t.Cell(13, 5).Range.Text = "96÷5=19, 1"

$t.Cell(17, 1).Range.Text = "18÷6=3, 0"
$t.Cell(17, 2).Range.Text = "22÷7=3, 1"
$t.Cell(17, 3).Range.Text = "68÷5=13, 3"
$t.Cell(17, 4).Range.Text = "32÷4=8, 0"
$t.Cell(17, 5).Range.Text = "92÷2=46, 0"
